# 自动更新价格数据: insert the newest day's reading as a new row 2 (just
# below the header), pushing the existing history down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (old row 2.."2025-12-06".. down to row 17
# "2025-11-21") down by one to make room for the new reading.
$ws.Rows(2).Insert()

# Force the date column to be read back as literal text (matching the
# rest of the column) instead of being auto-detected as a date serial.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-07"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Drop the temporary text format so the new row keeps the same (default)
# styling as every other data row instead of inheriting the header's bold
# border formatting that Insert() copied down.
$ws.Range("A2:D2").ClearFormats()
